{"js": "// Replace the two-digit multiplication problems in the practice table.\n// Each old problem string is unique in the document, so a plain\n// search-and-replace (matchCase, whole string incl. the trailing \"=\")\n// unambiguously targets the correct cell.\nconst replacements = [\n  [\"81\u00d748=\", \"30\u00d747=\"],\n  [\"48\u00d714=\", \"32\u00d762=\"],\n  [\"98\u00d727=\", \"75\u00d739=\"],\n  [\"21\u00d772=\", \"44\u00d789=\"],\n  [\"53\u00d767=\", \"51\u00d762=\"],\n  [\"28\u00d770=\", \"11\u00d715=\"],\n  [\"56\u00d750=\", \"14\u00d799=\"],\n  [\"63\u00d763=\", \"25\u00d785=\"],\n  [\"35\u00d789=\", \"11\u00d750=\"],\n  [\"27\u00d771=\", \"17\u00d768=\"],\n  [\"18\u00d741=\", \"66\u00d722=\"],\n  [\"54\u00d738=\", \"53\u00d796=\"],\n  [\"26\u00d736=\", \"47\u00d741=\"],\n  [\"79\u00d724=\", \"14\u00d762=\"],\n  [\"48\u00d755=\", \"88\u00d756=\"],\n  [\"55\u00d794=\", \"56\u00d799=\"],\n  [\"75\u00d764=\", \"72\u00d716=\"],\n  [\"57\u00d725=\", \"49\u00d752=\"],\n  [\"89\u00d786=\", \"11\u00d754=\"],\n  [\"31\u00d747=\", \"34\u00d744=\"],\n  [\"68\u00d719=\", \"37\u00d763=\"],\n  [\"74\u00d788=\", \"34\u00d763=\"],\n  [\"33\u00d738=\", \"42\u00d726=\"],\n  [\"17\u00d747=\", \"49\u00d771=\"],\n  [\"94\u00d784=\", \"81\u00d795=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems in the practice table.\n# Each old problem string is unique in the document, so a plain\n# Find/Replace on the whole \"NN\u00d7NN=\" string unambiguously targets the\n# correct cell.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"81\u00d748=\", \"30\u00d747=\"),\n    @(\"48\u00d714=\", \"32\u00d762=\"),\n    @(\"98\u00d727=\", \"75\u00d739=\"),\n    @(\"21\u00d772=\", \"44\u00d789=\"),\n    @(\"53\u00d767=\", \"51\u00d762=\"),\n    @(\"28\u00d770=\", \"11\u00d715=\"),\n    @(\"56\u00d750=\", \"14\u00d799=\"),\n    @(\"63\u00d763=\", \"25\u00d785=\"),\n    @(\"35\u00d789=\", \"11\u00d750=\"),\n    @(\"27\u00d771=\", \"17\u00d768=\"),\n    @(\"18\u00d741=\", \"66\u00d722=\"),\n    @(\"54\u00d738=\", \"53\u00d796=\"),\n    @(\"26\u00d736=\", \"47\u00d741=\"),\n    @(\"79\u00d724=\", \"14\u00d762=\"),\n    @(\"48\u00d755=\", \"88\u00d756=\"),\n    @(\"55\u00d794=\", \"56\u00d799=\"),\n    @(\"75\u00d764=\", \"72\u00d716=\"),\n    @(\"57\u00d725=\", \"49\u00d752=\"),\n    @(\"89\u00d786=\", \"11\u00d754=\"),\n    @(\"31\u00d747=\", \"34\u00d744=\"),\n    @(\"68\u00d719=\", \"37\u00d763=\"),\n    @(\"74\u00d788=\", \"34\u00d763=\"),\n    @(\"33\u00d738=\", \"42\u00d726=\"),\n    @(\"17\u00d747=\", \"49\u00d771=\"),\n    @(\"94\u00d784=\", \"81\u00d795=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
